$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")
$ws.Range("A24").NumberFormat = "@"
Write-Output $ws.Range("A24").NumberFormat
$ws.Range("A24:C24").Borders.Item(1).LineStyle = 1
$ws.Range("A24:C24").Borders.Item(1).Color = 0
Write-Output "ok"
